# mosip_master/xlsx/device_spec.xlsx
# "sierra leone master data"
#
# The device_spec master sheet was re-purposed from Madagascar to Sierra
# Leone: the localisation column (lang_code) switches from the French
# locale ("fra") to the English locale ("eng") used in Sierra Leone, and
# the data block gets an AutoFilter turned on (with its accompanying
# hidden _FilterDatabase defined name) over the sheet's used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount  = $usedRange.Rows.Count
$colCount  = $usedRange.Columns.Count

# Locate the lang_code column from row 1 headers (defensive: the sheet is
# A1:I6 today with lang_code in column A, but key off the header text).
$langCol = 1
for ($c = 1; $c -le $colCount; $c++) {
    if ($ws.Cells.Item(1, $c).Value2 -eq "lang_code") {
        $langCol = $c
        break
    }
}

# Every data row's language code moves from French to English.
for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, $langCol)
    if ($cell.Value2 -eq "fra") {
        $cell.Value2 = "eng"
    }
}

# Turn on the AutoFilter over the full data range (A1:I6) and make sure
# the workbook carries the usual hidden _xlnm._FilterDatabase name that
# Excel writes alongside an active AutoFilter.
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($rowCount, $colCount))
$dataRange.AutoFilter()

$filterDbRef  = $ws.Name + "!" + $dataRange.Address($true, $true)
$filterDbName = $ws.Names.Add("_xlnm._FilterDatabase", "=" + $filterDbRef)
$filterDbName.Visible = $false
